$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '38.808.14'
$ws.Range('E2').Value = '  +1.63%  '
$ws.Range('D3').Value = '2.095.08'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '229.18'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.34%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.613'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.29%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '61.41'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.75%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.387'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.67%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0846'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.60%  '
$ws.Range('E11').Value = '  -0.19%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '15.36'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +4.66%  '
$ws.Range('D13').Value = '2.404.75'
$ws.Range('E13').Value = '  -0.21%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '22.10'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.55%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.806'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +3.87%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.51'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.14%  '
$ws.Range('D17').Value = '2.082.06'
$ws.Range('E17').Value = '  -0.66%  '
$ws.Range('D18').Value = '38.739.11'
$ws.Range('E18').Value = '  +1.52%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '71.88'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +2.32%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.09'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.39%  '
$ws.Range('D21').Value = '0.0₃0843'
$ws.Range('E21').Value = '  +0.73%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '228.06'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.63%  '
$ws.Range('E23').Value = '  -0.49%  '
$ws.Range('E24').Value = '  -2.68%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.35'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.90%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '171.63'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.87%  '
$ws.Range('E27').Value = '  +1.02%  '
$ws.Range('E28').Value = '  +4.65%  '
$ws.Range('E29').Value = '  +5.87%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '19.35'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.47%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.45'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +2.69%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.121'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.48%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.53'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +2.16%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.76'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +1.54%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0614'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +1.22%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '6.48'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.27%  '
$ws.Range('E37').Value = '  -0.33%  '
$ws.Range('E38').Value = '  +1.48%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.999'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.16%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '18.03'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.16%  '
$ws.Range('E41').Value = '  +4.47%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '101.12'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +1.01%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '1.536.43'
$ws.Range('E43').Value = '  -0.71%  '
$ws.Range('E44').Value = '  -1.00%  '
$ws.Range('E45').Value = '  +0.59%  '
$ws.Range('B46').Value = 'TrustWalletToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.13'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.42%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '7.68'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +5.58%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '4.12'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.16%  '
$ws.Range('E49').Value = '  +1.25%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.97'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -1.14%  '
$ws.Range('D51').Value = '2.291.25'
$ws.Range('E51').Value = '  -0.18%  '
